# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer block that
# Jekyll's site build appended, along with the blank paragraph that preceded
# it. The final blank paragraph and the page-break paragraph are left intact.
$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $target1) {
        # the blank paragraph immediately preceding this one is also removed
        $startPara = $d.Paragraphs.Item($i - 1)
    }
    if ($text -eq $target2) {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
    Write-Output "Deleted footer block. Paragraphs remaining: $($d.Paragraphs.Count)"
} else {
    Write-Output "WARNING: footer block not found; no changes made."
}
